# Update "Inscritos" (E), "Pagos" (F), and "Inscrições homologadas" (H) figures
# in the Resumo de Inscricoes table to reflect the latest counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 66

$ws.Range("E5").Value = 70

$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 8
$ws.Range("H7").Value = 8

$ws.Range("E10").Value = 213
$ws.Range("F10").Value = 91
$ws.Range("H10").Value = 91

$ws.Range("E11").Value = 160

$ws.Range("E12").Value = 232
$ws.Range("F12").Value = 117
$ws.Range("H12").Value = 117

$ws.Range("E14").Value = 70

$ws.Range("E15").Value = 96

$ws.Range("E16").Value = 92
$ws.Range("F16").Value = 44
$ws.Range("H16").Value = 44

$ws.Range("E17").Value = 44

$ws.Range("E18").Value = 34

$ws.Range("E20").Value = 56

$ws.Range("E21").Value = 74

$ws.Range("E23").Value = 104

$ws.Range("E24").Value = 110

$ws.Range("E25").Value = 108
$ws.Range("F25").Value = 40
$ws.Range("H25").Value = 40

$ws.Range("E27").Value = 158
$ws.Range("F27").Value = 75
$ws.Range("H27").Value = 75

$ws.Range("E29").Value = 102

$ws.Range("E30").Value = 112
$ws.Range("F30").Value = 53
$ws.Range("H30").Value = 53

$ws.Range("E31").Value = 44
$ws.Range("F31").Value = 20
$ws.Range("H31").Value = 20

$ws.Range("E32").Value = 105
$ws.Range("F32").Value = 50
$ws.Range("H32").Value = 50

$ws.Range("E33").Value = 142

$ws.Range("E34").Value = 112
$ws.Range("F34").Value = 55
$ws.Range("H34").Value = 55

$ws.Range("E35").Value = 71
$ws.Range("F35").Value = 35
$ws.Range("H35").Value = 35

$ws.Range("E37").Value = 74
$ws.Range("F37").Value = 32
$ws.Range("H37").Value = 32

$ws.Range("E38").Value = 51
$ws.Range("F38").Value = 27
$ws.Range("H38").Value = 27

$ws.Range("E39").Value = 114
$ws.Range("F39").Value = 41
$ws.Range("H39").Value = 41

$ws.Range("E40").Value = 143

$ws.Range("E41").Value = 192
$ws.Range("F41").Value = 68
$ws.Range("H41").Value = 68

$ws.Range("E42").Value = 171
$ws.Range("F42").Value = 80
$ws.Range("H42").Value = 80

$ws.Range("E43").Value = 53

$ws.Range("E44").Value = 146

$ws.Range("E45").Value = 64
$ws.Range("F45").Value = 32
$ws.Range("H45").Value = 32

$ws.Range("E46").Value = 133
$ws.Range("F46").Value = 58
$ws.Range("H46").Value = 58

$ws.Range("E47").Value = 223

$ws.Range("E48").Value = 112

$ws.Range("E49").Value = 120

$ws.Range("E50").Value = 101

$ws.Range("E51").Value = 107
